$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 125, shifting the existing rows 125:153 down to 126:154
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly price observation
$ws.Cells.Item(125, 1).Value = 3
$ws.Cells.Item(125, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = 44476
$ws.Cells.Item(125, 5).Value = 5
$ws.Cells.Item(125, 6).Value = 100112001
$ws.Cells.Item(125, 7).Value = "Berenjena"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 60
$ws.Cells.Item(125, 11).Value = 9000
$ws.Cells.Item(125, 12).Value = 9000
$ws.Cells.Item(125, 13).Value = 9000
$ws.Cells.Item(125, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(125, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(125, 16).Value = 150
$ws.Cells.Item(125, 17).Value = 60
$ws.Cells.Item(125, 18).Value = "Hortaliza"
